{"js": "// Replace the two-digit-number / one-digit-number division prompts in the\n// worksheet table with a new set of problems (regenerated master data).\nconst replacements = [\n  [\"66\u00f77=\", \"10\u00f77=\"],\n  [\"43\u00f73=\", \"83\u00f77=\"],\n  [\"55\u00f78=\", \"80\u00f77=\"],\n  [\"93\u00f76=\", \"44\u00f78=\"],\n  [\"43\u00f78=\", \"48\u00f73=\"],\n  [\"11\u00f76=\", \"36\u00f77=\"],\n  [\"10\u00f74=\", \"82\u00f78=\"],\n  [\"41\u00f78=\", \"89\u00f76=\"],\n  [\"76\u00f77=\", \"30\u00f74=\"],\n  [\"67\u00f76=\", \"57\u00f77=\"],\n  [\"42\u00f78=\", \"78\u00f79=\"],\n  [\"79\u00f78=\", \"80\u00f75=\"],\n  [\"63\u00f75=\", \"65\u00f77=\"],\n  [\"78\u00f76=\", \"54\u00f76=\"],\n  [\"73\u00f79=\", \"59\u00f73=\"],\n  [\"57\u00f73=\", \"10\u00f77=\"],\n  [\"39\u00f74=\", \"81\u00f78=\"],\n  [\"61\u00f78=\", \"60\u00f77=\"],\n  [\"12\u00f76=\", \"52\u00f75=\"],\n  [\"54\u00f78=\", \"98\u00f79=\"],\n  [\"68\u00f78=\", \"52\u00f79=\"],\n  [\"85\u00f72=\", \"95\u00f78=\"],\n  [\"50\u00f75=\", \"51\u00f75=\"],\n  [\"36\u00f78=\", \"27\u00f75=\"],\n  [\"41\u00f73=\", \"68\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n  // Each search string is unique in this document, so replace the first hit.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-number / one-digit-number division prompts in the\n# worksheet table with a new set of problems (regenerated master data).\n$replacements = @(\n    @{ Old = \"66\u00f77=\"; New = \"10\u00f77=\" },\n    @{ Old = \"43\u00f73=\"; New = \"83\u00f77=\" },\n    @{ Old = \"55\u00f78=\"; New = \"80\u00f77=\" },\n    @{ Old = \"93\u00f76=\"; New = \"44\u00f78=\" },\n    @{ Old = \"43\u00f78=\"; New = \"48\u00f73=\" },\n    @{ Old = \"11\u00f76=\"; New = \"36\u00f77=\" },\n    @{ Old = \"10\u00f74=\"; New = \"82\u00f78=\" },\n    @{ Old = \"41\u00f78=\"; New = \"89\u00f76=\" },\n    @{ Old = \"76\u00f77=\"; New = \"30\u00f74=\" },\n    @{ Old = \"67\u00f76=\"; New = \"57\u00f77=\" },\n    @{ Old = \"42\u00f78=\"; New = \"78\u00f79=\" },\n    @{ Old = \"79\u00f78=\"; New = \"80\u00f75=\" },\n    @{ Old = \"63\u00f75=\"; New = \"65\u00f77=\" },\n    @{ Old = \"78\u00f76=\"; New = \"54\u00f76=\" },\n    @{ Old = \"73\u00f79=\"; New = \"59\u00f73=\" },\n    @{ Old = \"57\u00f73=\"; New = \"10\u00f77=\" },\n    @{ Old = \"39\u00f74=\"; New = \"81\u00f78=\" },\n    @{ Old = \"61\u00f78=\"; New = \"60\u00f77=\" },\n    @{ Old = \"12\u00f76=\"; New = \"52\u00f75=\" },\n    @{ Old = \"54\u00f78=\"; New = \"98\u00f79=\" },\n    @{ Old = \"68\u00f78=\"; New = \"52\u00f79=\" },\n    @{ Old = \"85\u00f72=\"; New = \"95\u00f78=\" },\n    @{ Old = \"50\u00f75=\"; New = \"51\u00f75=\" },\n    @{ Old = \"36\u00f78=\"; New = \"27\u00f75=\" },\n    @{ Old = \"41\u00f73=\"; New = \"68\u00f79=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($pair.Old)\"\n    }\n}\n"}
